$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1770.5
$ws.Range("I40").Value = 1400
$ws.Range("J40").Value = 2326.25
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 2326.25
$ws.Range("M40").Value = -1225
$ws.Range("N40").Value = -2676.25

# Row 41
$ws.Range("H41").Value = 7407767
$ws.Range("I41").Value = 11111392
$ws.Range("J41").Value = 516
$ws.Range("K41").Value = 11111392
$ws.Range("L41").Value = 516
$ws.Range("M41").Value = -11110952
$ws.Range("N41").Value = -1396

# Row 135
$ws.Range("H135").Value = 5925.913
$ws.Range("I135").Value = 6465.85
$ws.Range("J135").Value = 2326.3333
$ws.Range("K135").Value = 58192.65
$ws.Range("L135").Value = 20936.9997
$ws.Range("M135").Value = -55657.65
$ws.Range("N135").Value = -26006.9997

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2675.8948
$ws.Range("I61").Value = 1887.4
$ws.Range("J61").Value = 5632.75
$ws.Range("K61").Value = 1887.4
$ws.Range("L61").Value = 5632.75
$ws.Range("M61").Value = -1675.4
$ws.Range("N61").Value = -6056.75

# Row 136
$ws.Range("H136").Value = 2675.8948
$ws.Range("I136").Value = 1887.4
$ws.Range("J136").Value = 5632.75
$ws.Range("K136").Value = 5662.200000000001
$ws.Range("L136").Value = 16898.25
$ws.Range("M136").Value = -3112.200000000001
$ws.Range("N136").Value = -21998.25

$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 56071
$ws.Range("J20").Value = 56071
$ws.Range("L20").Value = 56071
$ws.Range("N20").Value = -56543

# Row 30
$ws.Range("H30").Value = 56071
$ws.Range("J30").Value = 56071
$ws.Range("L30").Value = 56071
$ws.Range("N30").Value = -56253

# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# Row 128
$ws.Range("H128").Value = 56071
$ws.Range("J128").Value = 56071
$ws.Range("L128").Value = 56071
$ws.Range("N128").Value = -66031

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

# Row 68
$ws.Range("H68").Value = 4190.39
$ws.Range("J68").Value = 2109.3333
$ws.Range("L68").Value = 6327.999899999999
$ws.Range("N68").Value = -7949.999899999999

# Row 71
$ws.Range("H71").Value = 4190.39
$ws.Range("J71").Value = 2109.3333
$ws.Range("L71").Value = 18983.9997
$ws.Range("N71").Value = -27095.9997

# Row 75
$ws.Range("H75").Value = 2500
$ws.Range("J75").Value = 3000
$ws.Range("L75").Value = 9000
$ws.Range("N75").Value = -10996

# Row 76
$ws.Range("H76").Value = 3500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 10500
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -11266

# Row 78
$ws.Range("H78").Value = 2500
$ws.Range("J78").Value = 3000
$ws.Range("L78").Value = 27000
$ws.Range("N78").Value = -36984

# Row 79
$ws.Range("H79").Value = 3500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 10500
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -13152

# Row 80
$ws.Range("H80").Value = 1150
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1150
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3450
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5322

# Row 81
$ws.Range("H81").Value = 7118075
$ws.Range("I81").Value = 1200
$ws.Range("J81").Value = 8541450
$ws.Range("K81").Value = 3600
$ws.Range("L81").Value = 25624350
$ws.Range("M81").Value = -2477
$ws.Range("N81").Value = -25626596

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 83
$ws.Range("H83").Value = 1150
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1150
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 10350
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -19710

# Row 84
$ws.Range("H84").Value = 7118075
$ws.Range("I84").Value = 1200
$ws.Range("J84").Value = 8541450
$ws.Range("K84").Value = 10800
$ws.Range("L84").Value = 76873050
$ws.Range("M84").Value = -5184
$ws.Range("N84").Value = -76884282

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# Row 132
$ws.Range("H132").Value = 797.0417
$ws.Range("I132").Value = 567.4
$ws.Range("J132").Value = 961.0714
$ws.Range("K132").Value = 5106.599999999999
$ws.Range("L132").Value = 8649.642600000001
$ws.Range("M132").Value = -2576.599999999999
$ws.Range("N132").Value = -13709.6426

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1987.091
$ws.Range("I102").Value = 1546
$ws.Range("K102").Value = 1546
$ws.Range("M102").Value = 76

# Row 122
$ws.Range("H122").Value = 1853518.5
$ws.Range("I122").Value = 3705037
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 11115111
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -11112661
$ws.Range("N122").Value = -10900

# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 92007.82000000001
$ws.Range("I126").Value = 125636.375
$ws.Range("J126").Value = 2331.6667
$ws.Range("K126").Value = 376909.125
$ws.Range("L126").Value = 6995.000100000001
$ws.Range("M126").Value = -374439.125
$ws.Range("N126").Value = -11935.0001

# Row 132
$ws.Range("H132").Value = 14707827
$ws.Range("I132").Value = 20001104
$ws.Range("J132").Value = 4282.222
$ws.Range("K132").Value = 60003312
$ws.Range("L132").Value = 12846.666
$ws.Range("M132").Value = -60000782
$ws.Range("N132").Value = -17906.666

Write-Output "Edits applied"